# Auto-generated edit script: refresh market-price derived columns (H-N)
# for the Leve profit rows, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3089966.8
$ws.Range("J17").Value = 3337092.2
$ws.Range("L17").Value = 10011276.6
$ws.Range("N17").Value = -10011612.6

$ws.Range("H64").Value = 3041.318
$ws.Range("I64").Value = 2980
$ws.Range("J64").Value = 3076.3572
$ws.Range("K64").Value = 2980
$ws.Range("L64").Value = 3076.3572
$ws.Range("M64").Value = -2732
$ws.Range("N64").Value = -3572.3572

$ws.Range("H67").Value = 3041.318
$ws.Range("I67").Value = 2980
$ws.Range("J67").Value = 3076.3572
$ws.Range("K67").Value = 2980
$ws.Range("L67").Value = 3076.3572
$ws.Range("M67").Value = -2122
$ws.Range("N67").Value = -4792.3572

$ws.Range("H76").Value = 3127.4375
$ws.Range("I76").Value = 2937.5833
$ws.Range("K76").Value = 2937.5833
$ws.Range("M76").Value = -2622.5833

$ws.Range("H79").Value = 3127.4375
$ws.Range("I79").Value = 2937.5833
$ws.Range("K79").Value = 2937.5833
$ws.Range("M79").Value = -1845.5833

$ws.Range("H112").Value = 940.95
$ws.Range("J112").Value = 1007
$ws.Range("L112").Value = 3021
$ws.Range("N112").Value = -5237

$ws.Range("H138").Value = 3071
$ws.Range("I138").Value = 1515.7727
$ws.Range("J138").Value = 3704.611
$ws.Range("K138").Value = 4547.3181
$ws.Range("L138").Value = 11113.833
$ws.Range("M138").Value = 592.6818999999996
$ws.Range("N138").Value = -21393.833

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4068.8667
$ws.Range("I122").Value = 3862.0833
$ws.Range("J122").Value = 4896
$ws.Range("K122").Value = 11586.2499
$ws.Range("L122").Value = 14688
$ws.Range("M122").Value = -9136.249899999999
$ws.Range("N122").Value = -19588

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1880.8448
$ws.Range("I31").Value = 1236.2
$ws.Range("J31").Value = 2571.5356
$ws.Range("K31").Value = 1236.2
$ws.Range("L31").Value = 2571.5356
$ws.Range("M31").Value = -941.2
$ws.Range("N31").Value = -3161.5356

$ws.Range("H34").Value = 1880.8448
$ws.Range("I34").Value = 1236.2
$ws.Range("J34").Value = 2571.5356
$ws.Range("K34").Value = 1236.2
$ws.Range("L34").Value = 2571.5356
$ws.Range("M34").Value = -1034.2
$ws.Range("N34").Value = -2975.5356

$ws.Range("H62").Value = 3275.8462
$ws.Range("I62").Value = 2647.5
$ws.Range("J62").Value = 3555.111
$ws.Range("K62").Value = 2647.5
$ws.Range("L62").Value = 3555.111
$ws.Range("M62").Value = -2023.5
$ws.Range("N62").Value = -4803.111

$ws.Range("H65").Value = 3275.8462
$ws.Range("I65").Value = 2647.5
$ws.Range("J65").Value = 3555.111
$ws.Range("K65").Value = 13237.5
$ws.Range("L65").Value = 17775.555
$ws.Range("M65").Value = -10117.5
$ws.Range("N65").Value = -24015.555

$ws.Range("H99").Value = 3514.7
$ws.Range("I99").Value = 3887.5386
$ws.Range("J99").Value = 2822.2856
$ws.Range("K99").Value = 3887.5386
$ws.Range("L99").Value = 2822.2856
$ws.Range("M99").Value = -2389.5386
$ws.Range("N99").Value = -5818.2856

$ws.Range("H107").Value = 1122.8889
$ws.Range("I107").Value = 1444.5
$ws.Range("J107").Value = 865.6
$ws.Range("K107").Value = 1444.5
$ws.Range("L107").Value = 865.6
$ws.Range("M107").Value = 475.5
$ws.Range("N107").Value = -4705.6

$ws.Range("H126").Value = 3514.7
$ws.Range("I126").Value = 3887.5386
$ws.Range("J126").Value = 2822.2856
$ws.Range("K126").Value = 11662.6158
$ws.Range("L126").Value = 8466.856800000001
$ws.Range("M126").Value = -9192.6158
$ws.Range("N126").Value = -13406.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4376.1035
$ws.Range("I70").Value = 2443.8572
$ws.Range("J70").Value = 4990.909
$ws.Range("K70").Value = 7331.571599999999
$ws.Range("L70").Value = 14972.727
$ws.Range("M70").Value = -7016.571599999999
$ws.Range("N70").Value = -15602.727

$ws.Range("H73").Value = 4376.1035
$ws.Range("I73").Value = 2443.8572
$ws.Range("J73").Value = 4990.909
$ws.Range("K73").Value = 7331.571599999999
$ws.Range("L73").Value = 14972.727
$ws.Range("M73").Value = -6239.571599999999
$ws.Range("N73").Value = -17156.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 92161830
$ws.Range("I70").Value = 414705900
$ws.Range("J70").Value = 6385.7144
$ws.Range("K70").Value = 414705900
$ws.Range("L70").Value = 6385.7144
$ws.Range("M70").Value = -414705630
$ws.Range("N70").Value = -6925.7144

$ws.Range("H73").Value = 92161830
$ws.Range("I73").Value = 414705900
$ws.Range("J73").Value = 6385.7144
$ws.Range("K73").Value = 414705900
$ws.Range("L73").Value = 6385.7144
$ws.Range("M73").Value = -414704964
$ws.Range("N73").Value = -8257.714400000001

$ws.Range("H80").Value = 3193.4482
$ws.Range("I80").Value = 4467.778
$ws.Range("J80").Value = 2620
$ws.Range("K80").Value = 4467.778
$ws.Range("L80").Value = 2620
$ws.Range("M80").Value = -3469.778
$ws.Range("N80").Value = -4616

$ws.Range("H83").Value = 3193.4482
$ws.Range("I83").Value = 4467.778
$ws.Range("J83").Value = 2620
$ws.Range("K83").Value = 22338.89
$ws.Range("L83").Value = 13100
$ws.Range("M83").Value = -17346.89
$ws.Range("N83").Value = -23084

$ws.Range("H102").Value = 1001.4048
$ws.Range("I102").Value = 881.1818
$ws.Range("J102").Value = 1442.2222
$ws.Range("K102").Value = 881.1818
$ws.Range("L102").Value = 1442.2222
$ws.Range("M102").Value = 740.8182
$ws.Range("N102").Value = -4686.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1893.45
$ws.Range("I7").Value = 1481.7273
$ws.Range("J7").Value = 2396.6667
$ws.Range("K7").Value = 1481.7273
$ws.Range("L7").Value = 2396.6667
$ws.Range("M7").Value = -1369.7273
$ws.Range("N7").Value = -2620.6667

$ws.Range("H40").Value = 2296
$ws.Range("I40").Value = 2296
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2296
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2160
$ws.Range("N40").ClearContents()

$ws.Range("H82").Value = 4330785.5
$ws.Range("I82").Value = 10101944
$ws.Range("J82").Value = 2416.6667
$ws.Range("K82").Value = 10101944
$ws.Range("L82").Value = 2416.6667
$ws.Range("M82").Value = -10101583
$ws.Range("N82").Value = -3138.6667

$ws.Range("H85").Value = 4330785.5
$ws.Range("I85").Value = 10101944
$ws.Range("J85").Value = 2416.6667
$ws.Range("K85").Value = 10101944
$ws.Range("L85").Value = 2416.6667
$ws.Range("M85").Value = -10100696
$ws.Range("N85").Value = -4912.6667

$ws.Range("H100").Value = 1718.75
$ws.Range("I100").Value = 1430
$ws.Range("J100").Value = 2200
$ws.Range("K100").Value = 1430
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -889
$ws.Range("N100").Value = -3282

$ws.Range("H126").Value = 1893.45
$ws.Range("I126").Value = 1481.7273
$ws.Range("J126").Value = 2396.6667
$ws.Range("K126").Value = 4445.1819
$ws.Range("L126").Value = 7190.000100000001
$ws.Range("M126").Value = -1975.1819
$ws.Range("N126").Value = -12130.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1443.8096
$ws.Range("I107").Value = 1510.4445
$ws.Range("J107").Value = 1393.8334
$ws.Range("K107").Value = 4531.333500000001
$ws.Range("L107").Value = 4181.5002
$ws.Range("M107").Value = -2611.333500000001
$ws.Range("N107").Value = -8021.5002
